$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44707
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 12500
$ws.Range("Q2").Value = '$/caja 12 kilos empedrada'
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1042
$ws.Range("T2").Value = 12

# Row 3
$ws.Range("D3").Value = 44708
$ws.Range("M3").Value = 70
$ws.Range("P3").Value = 12571
$ws.Range("S3").Value = 1048

# Row 4
$ws.Range("D4").Value = 44742
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 806
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44334
$ws.Range("L5").Value = "Primera"
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 11500
$ws.Range("Q5").Value = '$/caja 12 kilos granel'
$ws.Range("S5").Value = 11500
$ws.Range("T5").Value = 1

# Row 6
$ws.Range("D6").Value = 44714
$ws.Range("M6").Value = 100
$ws.Range("P6").Value = 14500
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 806

# Row 8
$ws.Range("D8").Value = 44719
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14400
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = "Región del Maule"
$ws.Range("S8").Value = 800
$ws.Range("T8").Value = 18
